$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- Row 1: updated date serial --
$ws.Range("A1").Value = 45413

# -- Column F labels (right-aligned), rows 2-23 --
$ws.Range("F2").Value = "RAW MATERIALS:"
$ws.Range("F2").HorizontalAlignment = -4152

$ws.Range("F3").HorizontalAlignment = -4152

$ws.Range("F4").Value = "HDPE Granules"
$ws.Range("F4").HorizontalAlignment = -4152

$ws.Range("F5").Value = "Master Batches"
$ws.Range("F5").HorizontalAlignment = -4152

$ws.Range("F6").Value = "Colour Pigments"
$ws.Range("F6").HorizontalAlignment = -4152

$ws.Range("F7").Value = "Total"
$ws.Range("F7").HorizontalAlignment = -4152

$ws.Range("F8").Value = "Work in Progress"
$ws.Range("F8").HorizontalAlignment = -4152

$ws.Range("F9").Value = "HDPE Tape - Factory"
$ws.Range("F9").HorizontalAlignment = -4152

$ws.Range("F10").Value = "HDPE Tape - Job Work"
$ws.Range("F10").HorizontalAlignment = -4152

$ws.Range("F11").Value = "Total"
$ws.Range("F11").HorizontalAlignment = -4152

$ws.Range("F12").Value = "FINISHED GOODS:"
$ws.Range("F12").HorizontalAlignment = -4152

$ws.Range("F13").Value = "HDPE Fishnet Fabrics"
$ws.Range("F13").HorizontalAlignment = -4152

$ws.Range("F14").Value = "Shadenet Fabrics & Weed Mat"
$ws.Range("F14").HorizontalAlignment = -4152

$ws.Range("F15").Value = "PP Fabric & Sacks"
$ws.Range("F15").HorizontalAlignment = -4152

$ws.Range("F16").Value = "Total"
$ws.Range("F16").HorizontalAlignment = -4152

$ws.Range("F17").Value = "Consumbles and Thread"
$ws.Range("F17").HorizontalAlignment = -4152

$ws.Range("F18").Value = "Packing Materials"
$ws.Range("F18").HorizontalAlignment = -4152

$ws.Range("F19").Value = "Seconds"
$ws.Range("F19").HorizontalAlignment = -4152

$ws.Range("F20").Value = "Total"
$ws.Range("F20").HorizontalAlignment = -4152

$ws.Range("F21").Value = "Grand Total"
$ws.Range("F21").HorizontalAlignment = -4152

$ws.Range("F22").Value = "SFG/FG"
$ws.Range("F22").HorizontalAlignment = -4152

$ws.Range("F23").Value = "Diff in SFG/FG"
$ws.Range("F23").HorizontalAlignment = -4152

# -- Raw materials section --
$ws.Range("A4").Value = 191250
$ws.Range("B4").Value = 98.297767320261443
$ws.Range("C4").Value = 18799448

$ws.Range("A5").Value = 2626.0000000000005
$ws.Range("B5").Value = 217.47220106626042
$ws.Range("C5").Value = 571082

$ws.Range("B6").Value = 2243.6888888888889
$ws.Range("C6").Value = 100966

$ws.Range("A7").Value = 193921
$ws.Range("B7").Value = 100.40942445635078
$ws.Range("C7").Value = 19471496

# -- Work in progress section --
$ws.Range("A9").Value = 47620
$ws.Range("B9").Value = 163
$ws.Range("C9").Value = 7762060

$ws.Range("A10").Value = 57609
$ws.Range("B10").Value = 163
$ws.Range("C10").Value = 9390267

$ws.Range("A11").Value = 105229
$ws.Range("B11").ClearContents()
$ws.Range("C11").Value = 17152311

# -- Finished goods section --
$ws.Range("A13").Value = 263759
$ws.Range("B13").Value = 219
$ws.Range("C13").Value = 57763289

$ws.Range("A14").Value = 1056
$ws.Range("B14").Value = 200.96022727272728
$ws.Range("C14").Value = 212214

$ws.Range("A16").Value = 264815
$ws.Range("B16").ClearContents()
$ws.Range("C16").Value = 57975503

# -- Grand total section --
$ws.Range("A21").Value = 563965
$ws.Range("B21").Value = 167.73968242710097
$ws.Range("C21").Value = 94599310

$ws.Range("A22").Value = 370044
$ws.Range("B22").Value = 203.02400255104797
$ws.Range("C22").Value = 75127814

$ws.Range("A23").Value = 24642
$ws.Range("C23").Value = 8136504

# -- Selection --
$ws.Range("A2").Select()
